$d = $word.ActiveDocument

# Third table in the document: the "signature" table listing the manager
# and job-title placeholders. It currently auto-fits its width; switch it
# to a fixed width and widen its third (rightmost) column.
$tbl = $d.Tables(3)

# Widen the third column from 3645 dxa (182.25pt) to 3933 dxa (196.65pt).
# Word propagates a cell-width write to every cell in that column (both
# rows here) and updates the corresponding <w:gridCol/> entry too.
$tbl.Cell(1, 3).Width = 196.65

# Switch the table from automatic width to a fixed width equal to the sum
# of the (now updated) column widths: 425 + 4293 + 3933 = 8651 dxa, i.e.
# 432.55 points.
$tbl.PreferredWidthType = 3
$tbl.PreferredWidth = 432.55
